$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A40").Value = "Record"
$ws.Range("B40").Value = "RJ No Ar TV Record"
$ws.Range("C40").Value = "Cultura"
$ws.Range("D40").Value = "2025-04-02T11:15"
$ws.Range("E40").Value = "Negativo"
$ws.Range("F40").Value = "testando22"
